$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Thin Cerablanket
$ws.Range("D2").Value = 278.41
$ws.Range("E2").Value = 171.59
$ws.Range("F2").Value = 40.37411764705882

# Row 3 - Medium Rock Wool
$ws.Range("D3").Value = 246.83
$ws.Range("E3").Value = 203.17
$ws.Range("F3").Value = 47.80470588235294

# Row 4 - Thick Silika Mat
$ws.Range("D4").Value = 238.83
$ws.Range("E4").Value = 211.17
$ws.Range("F4").Value = 49.68705882352941

# Row 5 - Multi-layer Needeled
$ws.Range("D5").Value = 187.41
$ws.Range("E5").Value = 262.59
$ws.Range("F5").Value = 61.78588235294118
